# Add a new weekly record at the top of the Higo data set (row 16),
# pushing all existing rows 16-32 down by one (new last row becomes 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 16, shifting rows 16:32 down to 17:33.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Cells.Item(16, 1).Value2 = 6
$ws.Cells.Item(16, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value2 = 44679
$ws.Cells.Item(16, 5).Value2 = 13
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value2 = 100101
$ws.Cells.Item(16, 8).Value = "Berries"
$ws.Cells.Item(16, 9).Value2 = 100101006
$ws.Cells.Item(16, 10).Value = "Higo"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value2 = 150
$ws.Cells.Item(16, 14).Value2 = 12000
$ws.Cells.Item(16, 15).Value2 = 12000
$ws.Cells.Item(16, 16).Value2 = 12000
$ws.Cells.Item(16, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(16, 18).Value = "Región Metropolitana"
$ws.Cells.Item(16, 19).Value2 = 1714
$ws.Cells.Item(16, 20).Value2 = 7
